$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New card-data rows appended below the existing A1:G7 data (sheet grows to A1:G10).
# Columns: A=Name, B=Email, C=CardNumber, D=ExpirationDate, E=SecurityCode,
#          F=Timestamp, G=Phone
# Several values look numeric (SecurityCode, Phone) but must stay text, matching
# how the rest of the sheet stores every value as text.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 8
$ws.Range("A8").Value = "WAssem"
$ws.Range("B8").Value = "was@jjj.com"
$ws.Range("C8").Value = "7007********5055"
$ws.Range("D8").Value = "12-∞"
Set-TextValue $ws.Range("E8") "985"
$ws.Range("F8").Value = "11/18/2025, 4:23:28 PM"
Set-TextValue $ws.Range("G8") "+963879394"

# Row 9
$ws.Range("A9").Value = "Ahmad Bassam Abboud"
$ws.Range("B9").Value = "student@university.edu"
$ws.Range("C9").Value = "7007********5055"
$ws.Range("D9").Value = "12-∞"
Set-TextValue $ws.Range("E9") "985"
$ws.Range("F9").Value = "11/18/2025, 4:23:58 PM"
Set-TextValue $ws.Range("G9") "+963879394"

# Row 10
$ws.Range("A10").Value = "test"
$ws.Range("B10").Value = "test@mail.com"
$ws.Range("C10").Value = "7007********5055"
$ws.Range("D10").Value = "12-∞"
Set-TextValue $ws.Range("E10") "985"
$ws.Range("F10").Value = "11/18/2025, 4:24:11 PM"
Set-TextValue $ws.Range("G10") "1723871893"

Write-Host "Added 3 rows of card data (rows 8-10)"
